# Update "想去人数" (F column) counts across the workbook's sheets.
# Mirrors the commit "Update gh-pages to output generated at 456a3b4":
# small counter increments on the F column of 展览, 演出, 本地生活 and 全部类型.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 351
$ws1.Range("F8").Value  = 66
$ws1.Range("F9").Value  = 546
$ws1.Range("F17").Value = 6728
$ws1.Range("F19").Value = 76
$ws1.Range("F20").Value = 23
$ws1.Range("F21").Value = 7647
$ws1.Range("F26").Value = 2152
$ws1.Range("F27").Value = 919
$ws1.Range("F35").Value = 1779
$ws1.Range("F39").Value = 7
$ws1.Range("F41").Value = 1249
$ws1.Range("F42").Value = 1881
$ws1.Range("F43").Value = 2152

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 73

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 245
$ws3.Range("F3").Value = 1242
$ws3.Range("F4").Value = 80

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 245
$ws4.Range("F4").Value  = 1242
$ws4.Range("F5").Value  = 80
$ws4.Range("F7").Value  = 351
$ws4.Range("F10").Value = 66
$ws4.Range("F11").Value = 546
$ws4.Range("F15").Value = 73
$ws4.Range("F20").Value = 6728
$ws4.Range("F22").Value = 76
$ws4.Range("F23").Value = 23
$ws4.Range("F24").Value = 7647
$ws4.Range("F29").Value = 2152
$ws4.Range("F30").Value = 919
$ws4.Range("F38").Value = 1779
$ws4.Range("F42").Value = 7
$ws4.Range("F44").Value = 1249
$ws4.Range("F45").Value = 1881
$ws4.Range("F47").Value = 2152
